$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need to be
# forced to Text format first, otherwise Excel auto-converts them to
# numeric cells (losing the original text-price formatting / trailing zeros).
$ws.Range('D2').Value = '28.107.44'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '1.789.54'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.33'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5173'
$ws.Range('E7').Value = '  +1.41%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3799'
$ws.Range('E8').Value = '  -3.79%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07994'
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.267'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.002'
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.45'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.269'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').Value = '1.788.01'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('E18').Value = '  -4.42%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06553'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.28'
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.944'
$ws.Range('E22').Value = '  -2.54%  '
$ws.Range('D23').Value = '28.146.02'
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.13'
$ws.Range('E24').Value = '  -3.06%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '160.49'
$ws.Range('E26').Value = '  +2.91%  '
$ws.Range('E27').Value = '  -4.24%  '
$ws.Range('D28').Value = '1.993.40'
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.324'
$ws.Range('E29').Value = '  -3.63%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '122.67'
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1081'
$ws.Range('E31').Value = '  -1.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.050'
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.666'
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.530'
$ws.Range('E34').Value = '  -4.41%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07173'
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.13'
$ws.Range('E36').Value = '  +7.34%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02306'
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2140'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.061'
$ws.Range('E39').Value = '  -3.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.555'
$ws.Range('E40').Value = '  -3.44%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6146'
$ws.Range('E41').Value = '  -2.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.160'
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.354'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.23'
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.762'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5931'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '127.49'
$ws.Range('E47').Value = '  +1.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.215'
$ws.Range('E48').Value = '  +2.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.915'
$ws.Range('E49').Value = '  -3.59%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06749'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '72.51'
$ws.Range('E51').Value = '  -2.48%  '
